$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date column (F) with header and values
$ws.Range("F1").Value = "25_03_2024"
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 4

# Match the underline styling applied to F2 and F5 (same as B5/D5)
$ws.Range("F2").Font.Underline = $true
$ws.Range("F5").Font.Underline = $true

# Update the active selection to F2, matching the post-edit state
$ws.Range("F2").Select()

# Touch page setup so a pageSetup element is emitted (portrait, as in the target)
$ws.PageSetup.Orientation = 1
